{"js": "const body = context.document.body;\n\n// 1) \"RESOLUCI\u00d3N \" -> \"RESOLUCI\u00d3N N\u00b0 \" (insert \"N\u00b0 \" right after the literal\n//    \"RESOLUCI\u00d3N \" label, before the RESOLUCION merge field).\nconst resolucionLabel = body.search(\"RESOLUCI\u00d3N \", { matchCase: true, matchWholeWord: false });\nresolucionLabel.load(\"text\");\nawait context.sync();\nresolucionLabel.items[0].insertText(\"N\u00b0 \", Word.InsertLocation.after);\n\n// 2) \"${RESOLUCION}\" -> \"${RESOLUCION_APROBACION}\" (cached mail-merge field\n//    result text; the underlying MERGEFIELD instruction is left untouched).\nconst resolucionField = body.search(\"${RESOLUCION}\", { matchCase: true, matchWholeWord: false });\nresolucionField.load(\"text\");\nawait context.sync();\nresolucionField.items[0].insertText(\"${RESOLUCION_APROBACION}\", Word.InsertLocation.replace);\n\n// 3) \"${FECHA_DE_RESOLUCION}\" -> \"${FECHA_DE_RESOLUCION_APROBACION}\"\nconst fechaField = body.search(\"${FECHA_DE_RESOLUCION}\", { matchCase: true, matchWholeWord: false });\nfechaField.load(\"text\");\nawait context.sync();\nfechaField.items[0].insertText(\"${FECHA_DE_RESOLUCION_APROBACION}\", Word.InsertLocation.replace);\n\n// 4) \"El expediente \" -> \"El expediente N\u00b0 \" (insert \"N\u00b0 \" before the\n//    N_MESA_DE_PARTES merge field).\nconst expedienteLabel = body.search(\"El expediente \", { matchCase: true, matchWholeWord: false });\nexpedienteLabel.load(\"text\");\nawait context.sync();\nexpedienteLabel.items[0].insertText(\"N\u00b0 \", Word.InsertLocation.after);\n\n// 5) \"Oficio \" -> \"Oficio N\u00b0 \" (insert \"N\u00b0 \" before the\n//    OFICIO_DE_PRESENTACION_FAC merge field), then relocate the \"_GoBack\"\n//    bookmark so it sits right after the newly inserted \"N\u00b0 \" text instead\n//    of at the very end of the document.\nconst oficioLabel = body.search(\"Oficio \", { matchCase: true, matchWholeWord: false });\noficioLabel.load(\"text\");\nawait context.sync();\nconst oficioNumRange = oficioLabel.items[0].insertText(\"N\u00b0 \", Word.InsertLocation.after);\nawait context.sync();\n\ncontext.document.deleteBookmark(\"_GoBack\");\nconst afterOficioNum = oficioNumRange.getRange(Word.RangeLocation.after);\nafterOficioNum.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"RESOLUCI\u00d3N \" -> \"RESOLUCI\u00d3N N\u00b0 \" (insert \"N\u00b0 \" right after the literal\n#    \"RESOLUCI\u00d3N \" label, before the RESOLUCION merge field).\n$range = $d.Content\n$find = $range.Find\n$find.MatchCase = $true\n$find.Text = \"RESOLUCI\u00d3N \"\n$find.Execute()\n$range.Collapse(0)\n$range.InsertAfter(\"N\u00b0 \")\n\n# 2) \"${RESOLUCION}\" -> \"${RESOLUCION_APROBACION}\" (cached mail-merge field\n#    result text; the underlying MERGEFIELD instruction is left untouched).\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.Text = '${RESOLUCION}'\n$find2.Replacement.Text = '${RESOLUCION_APROBACION}'\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n\n# 3) \"${FECHA_DE_RESOLUCION}\" -> \"${FECHA_DE_RESOLUCION_APROBACION}\"\n$range3 = $d.Content\n$find3 = $range3.Find\n$find3.Text = '${FECHA_DE_RESOLUCION}'\n$find3.Replacement.Text = '${FECHA_DE_RESOLUCION_APROBACION}'\n$find3.Execute($find3.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find3.Replacement.Text, 2)\n\n# 4) \"El expediente \" -> \"El expediente N\u00b0 \" (insert \"N\u00b0 \" before the\n#    N_MESA_DE_PARTES merge field).\n$range4 = $d.Content\n$find4 = $range4.Find\n$find4.MatchCase = $true\n$find4.Text = \"El expediente \"\n$find4.Execute()\n$range4.Collapse(0)\n$range4.InsertAfter(\"N\u00b0 \")\n\n# 5) \"Oficio \" -> \"Oficio N\u00b0 \" (insert \"N\u00b0 \" before the\n#    OFICIO_DE_PRESENTACION_FAC merge field), then relocate the \"_GoBack\"\n#    bookmark so it sits right after the newly inserted \"N\u00b0 \" text instead\n#    of at the very end of the document.\n$range5 = $d.Content\n$find5 = $range5.Find\n$find5.MatchCase = $true\n$find5.Text = \"Oficio \"\n$find5.Execute()\n$range5.Collapse(0)\n$range5.InsertAfter(\"N\u00b0 \")\n$range5.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $range5)\n"}
